# Applies the "update to share with Christopher" edit:
#  - bs_treatment_table (sheet2): rename headers (sample->run, treatment-> dropped,
#    description->id), move the old "description" values into a new "treatInfo"
#    column (B, styled w/ Menlo font), add a new "id" column (c1..c3,t1..t3) in C,
#    and a new numeric "treat" column (0/1) in D.
#  - rna_treatment_table (sheet3): rename headers only (sample->run,
#    treatment->treat, description->treatInfo); data values unchanged.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 2: bs_treatment_table
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("bs_treatment_table")

# Capture the old column B/C values (per-row) before overwriting them.
$oldTreatment = @{}
$oldDescription = @{}
for ($r = 2; $r -le 7; $r++) {
    $oldTreatment[$r] = $ws2.Cells.Item($r, 2).Value()
    $oldDescription[$r] = $ws2.Cells.Item($r, 3).Value()
}

# New header row. Order matters for shared-string allocation order:
# "treatInfo" must be interned before "run" to match upstream's table.
$ws2.Cells.Item(1, 2).Value = "treatInfo"
$ws2.Cells.Item(1, 1).Value = "run"
$ws2.Cells.Item(1, 3).Value = "id"
$ws2.Cells.Item(1, 4).Value = "treat"

# Style the new "treatInfo" header cell with the Menlo font used upstream.
$f = $ws2.Cells.Item(1, 2).Font
$f.Name = "Menlo"
$f.Size = 11
$f.Color = 0

$ids = @{2="c1"; 3="c2"; 4="c3"; 5="t1"; 6="t2"; 7="t3"}
$treatFlag = @{2=0; 3=0; 4=0; 5=1; 6=1; 7=1}

for ($r = 2; $r -le 7; $r++) {
    # Column B becomes the old "description" text (treatInfo).
    $ws2.Cells.Item($r, 2).Value = $oldDescription[$r]
    # Column C becomes the new short id.
    $ws2.Cells.Item($r, 3).Value = $ids[$r]
    # Column D is the new numeric treat flag.
    $ws2.Cells.Item($r, 4).Value = $treatFlag[$r]
}

$ws2.Range("A1:D7").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: rna_treatment_table
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("rna_treatment_table")

$ws3.Cells.Item(1, 1).Value = "run"
$ws3.Cells.Item(1, 2).Value = "treat"
$ws3.Cells.Item(1, 3).Value = "treatInfo"

$ws3.Range("A2").Select() | Out-Null
